$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.226.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.423.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.19%  "
$ws.Range("E4").Value = "  -0.73%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.427.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.79%  "
$ws.Range("E10").Value = "  +4.96%  "
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("E13").Value = "  +4.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.48%  "
$ws.Range("E15").Value = "  +10.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.862.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.841.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +22.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.431.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.74%  "
$ws.Range("E20").Value = "  +3.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.30%  "
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "562.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.46%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0936"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.38%  "
$ws.Range("E32").Value = "  +7.18%  "
$ws.Range("E33").Value = "  +4.52%  "
$ws.Range("E34").Value = "  +5.73%  "
$ws.Range("E35").Value = "  +4.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.50%  "
$ws.Range("E37").Value = "  +13.78%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  +5.50%  "
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "146.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "149.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.41%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.82%  "
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0542"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.68%  "
$ws.Range("E48").Value = "  +7.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.593"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0913"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("E51").Value = "  +3.98%  "
